# Add new localization entries for the photo album / photo result screens,
# plus the percent/points labels, as new rows appended to the Key/Value table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A108").Value = "photo_album"
$ws.Range("B108").Value = "Photo Album"

$ws.Range("A109").Value = "photo_result"
$ws.Range("B109").Value = "Photo Result"

# NOTE: "percent"/"Percent:" is entered before "points"/"Points:" even
# though it ends up one row below it, so that the shared-string table
# order matches how the workbook was actually authored.
$ws.Range("A111").Value = "percent"
$ws.Range("B111").Value = "Percent:"

$ws.Range("B110").Value = "Points:"
$ws.Range("A110").Value = "points"

$ws.Range("A111").Select()
